$d = $word.ActiveDocument

# Locate the final paragraph of the document - the one ending with
# "...control of master branch." which currently also holds the
# trailing _GoBack bookmark.
$target = $d.Content
$target.Find.Execute("control of master branch.", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Expand(4) | Out-Null

# Exclude the paragraph's own end-of-paragraph mark from the range so
# InsertXML replaces only the paragraph's contents (run + bookmark)
# without leaving a stray empty paragraph behind.
$target.MoveEnd(1, -1) | Out-Null

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00667208" w:rsidRDefault="00667208" w:rsidP="00B815B2"><w:r><w:t>In cases like that what we can do is git rebase&#8230;what it actually does that it takes up your complete branch which was to be merged to say master&#8230; and commits all its commits in master one after another so as we get one linear tree structure of it. But there is a drawback to it , that is it looses the author ship&#8230;all the commits would be under the name of the person who is having the control of master branch.</w:t></w:r></w:p>
<w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr></w:p>
<w:p><w:r><w:t>For admin n master:-</w:t></w:r></w:p>
<w:p><w:r><w:t>Git checkout admin</w:t></w:r></w:p>
<w:p><w:r><w:t>Git rebase master</w:t></w:r></w:p>
<w:p><w:r><w:t>Git checkout master</w:t></w:r></w:p>
<w:p><w:r><w:t>Git merge admin</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xml) | Out-Null
